$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old worked-example rows (5-12); rows 1-4 are kept and re-purposed
# below as the new, richer set of example prompts that get saved to github.
$ws.Rows("5:12").Delete()

# Row 2: "I can write a for loop" example, now wrapped with framing text for
# both the user example and the assistant's follow-up prompt.
$ws.Range("A2").Value = "Here is an example of a response from the user who has some knowledge in for loops as used in C.
 I can write a for loop"
$ws.Range("B2").Value = "Great to know that you have knowledge of the for loops. Lets put your knowledge into practice. Write a loop that will display  a range of numbers from 1 to 10"

# Row 3: "I don't know anything about loops" example, with new framing text.
$ws.Range("A3").Value = "Here is an example of a respose from the user with no knowledge of the for loops as used in C.
I don't know anything about loops"
$ws.Range("B3").Value = "It is okay that you do not know about loops. Let’s start by learning the concepts you already know and are  the building blocks of the for loop. We will start by  looking at concept of variables. Define a variable that will store the number 5?"

# Row 4: "loops help in repeating things in your code" example, with new framing text.
$ws.Range("A4").Value = "Here is a response from a user who has somewhat knowledge of for loops. 
loops help in repeating things in your code"

# B4 keeps its original text, unchanged.

# The example column now wraps its text like the assistant-prompt column does.
$ws.Range("A2:A4").WrapText = $true
